$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet and its matching defined name
# (old name "Recaptures___EDI_query" -> new name "qry_Knights_Recaptures_EDI")
$ws.Name = "qry_Knights_Recaptures_EDI"
$wb.Names.Item(1).Name = "qry_Knights_Recaptures_EDI"

# Set the width of column K (the 11th column) to a custom width.
# 10 + 2/7 character-width units serializes to the OOXML width="11".
$ws.Columns.Item(11).ColumnWidth = 10.2857142857142857

# Update F329 value from "Spring" to "Fall"
$ws.Range("F329").Value = "Fall"
